$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","D","E","F","G","H","I","L","N","O")

$data = @{
  2 = @(0.8293553357516998,0.008354403957902434,0.4009701891499269,0.368015581926592,0.2392252137583739,0.3487380636063904,0.6322407628076663,0.2442602145456902,1.369676114756047,1.089643772403861)
  3 = @(0.7728124514055992,0.007289734207173382,0.3996912454548962,0.3522789147792693,0.2241302352068573,0.3449989052252533,0.6433315749451842,0.2182181681758664,1.348640746459168,1.049535946167481)
  4 = @(0.7382075621112847,0.00663529901271076,0.3990165213267183,0.3428678498180844,0.2150318213095233,0.3429043927951767,0.6505462909663908,0.2022314481321672,1.336207966503693,1.025680892599183)
  5 = @(0.7241353756353988,0.006368444586602351,0.3987696282450486,0.3390959454812759,0.2113667816152827,0.3421015301690318,0.6535881450836105,0.1957179488867666,1.331264012739695,1.016153575886818)
  6 = @(0.7218005237320142,0.006324123985972818,0.3987303344786532,0.338473439225119,0.210760776822525,0.3419712768011749,0.6540993904589776,0.1946364731097532,1.330450502858355,1.014583273778726)
  7 = @(0.7380176582229865,0.006631700775479032,0.3990130776662575,0.3428167247985527,0.2149822207986745,0.3428933599060571,0.6505869023208088,0.2021435993192995,1.336140793203839,1.025551619331651)
  8 = @(0.8098366984607992,0.00798746484403523,0.4005064208445788,0.3625374056686539,0.2339851237067592,0.3474070332428738,0.6359808245350891,0.2352804926643586,1.362323484166254,1.075654407985724)
  9 = @(0.9515195880990177,0.01063982697426979,0.4043008682587654,0.4032067358948126,0.272605023676789,0.3578551815531483,0.6105516498545276,0.3002722091572707,1.41745822540139,1.18003880048937)
  10 = @(1.056075670166138,0.01258413275497361,0.4076021244115609,0.4343118377337305,0.3018175671259513,0.3665053041325308,0.5938278779896704,0.3480123772628474,1.460225878566689,1.26049619156052)
  11 = @(1.103730379419403,0.01346758797208736,0.4092124222222964,0.4487302905739341,0.3152919705290458,0.3706521012745014,0.5866452001819846,0.3697255831099255,1.480162495263897,1.297922490577264)
  12 = @(1.12178801706159,0.01380196940544209,0.4098375470151581,0.4542288718730418,0.3204211973422417,0.3722528273535204,0.5839864500823131,0.3779468679210822,1.487780240448217,1.312213901100478)
  13 = @(1.117898476734524,0.01372996194653808,0.4097022372134305,0.4530429360670922,0.3193153348558724,0.3719067297903962,0.5845563379181913,0.3761763204080921,1.486136606870019,1.309130700765508)
  14 = @(1.105215762649891,0.01349510111807461,0.4092635457107363,0.4491818880289458,0.3157134182227566,0.3707831844294418,0.5864252368906424,0.3704019765057467,1.480787851238546,1.299095870200546)
  15 = @(1.097448735368175,0.01335122022089053,0.4089968242790363,0.4468219149512294,0.3135106282651918,0.3700989415242759,0.5875779597049666,0.3668648758574022,1.477520432353174,1.292964732601888)
  16 = @(1.05296305175716,0.01252637489091768,0.4074990488351418,0.4333749629507082,0.3009407232612915,0.3662385604735761,0.5943058404501427,0.346593248112697,1.458932571346224,1.258066923962588)
  17 = @(1.025695036876243,0.01202008564928292,0.4066078208117432,0.4251944882691987,0.2932770636235489,0.36392456779474,0.598542080157257,0.3341559038144055,1.447652072999801,1.236869857913604)
  18 = @(1.010019899921645,0.01172878657389731,0.4061054422651864,0.4205145761161475,0.2888865813405062,0.3626135593426199,0.6010186577654917,0.3270019195981888,1.441209201860772,1.224755589655132)
  19 = @(1.004714103744334,0.01163014193991074,0.4059371110981544,0.4189343800700414,0.2874030340335025,0.3621730999800263,0.6018640536164854,0.3245796567386208,1.439035581205502,1.220667251585013)
  20 = @(1.028596872507876,0.01207399095006423,0.4067016365024756,0.4260626972898365,0.2940910663780585,0.3641688329080495,0.5980869845058123,0.3354799201476339,1.448848212938316,1.23911827672552)
  21 = @(1.108940672730057,0.01356409002243453,0.4093919858475061,0.4503149231113213,0.3167706613962906,0.3711123715158635,0.5858746352941928,0.3720980729937935,1.482357068278276,1.302040115720416)
  22 = @(1.161518344932631,0.01453699256988017,0.4112395740627832,0.4663902972515501,0.3317491221454958,0.3758276933411935,0.5782497529251724,0.39602390907055,1.504653975998735,1.343856173712624)
  23 = @(1.133450839858938,0.01401783005906765,0.4102453997441629,0.4577899663822791,0.3237405287952271,0.3732948242693794,0.5822866504599382,0.3832549615241874,1.492717719727921,1.321474705839449)
  24 = @(1.027284947505848,0.01204962106824325,0.4066591912595072,0.425670107876357,0.2937230075908559,0.364058340414033,0.5982926052097524,0.3348813435103182,1.448307305803496,1.238101541669522)
  25 = @(0.9131056473623858,0.009923022574177764,0.4031831563568247,0.3919899786681071,0.2620109472579486,0.3548576853212211,0.6170868968725438,0.2826907280841908,1.477520432353174,1.292964732601888)
}

foreach ($r in $data.Keys) {
  $rowVals = $data[$r]
  for ($i = 0; $i -lt $cols.Count; $i++) {
    $ws.Range("$($cols[$i])$r").Value = $rowVals[$i]
  }
}
